$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.001.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.568.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.51%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.23"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.567.06"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.29%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +1.25%  "

$ws.Range("E10").Value = "  +0.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.388"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.175.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000183"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.567.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.31%  "

$ws.Range("E17").Value = "  +0.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.717.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.581"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.712.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.62%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  +7.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.52%  "

$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.57%  "

$ws.Range("E32").Value = "  +22.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.576.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.51%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.145"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "169.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0810"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.826"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.492.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.15%  "
